# Apply the crypto price/rank refresh described in the commit:
# "Updated symbol list on Fri Dec 23 07:36:20 UTC 2022 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking Price values are written with a leading apostrophe so Excel
# stores them as literal text (matching the workbook's existing text-string cells)
# instead of silently converting them to numbers.

# Row 2
$ws.Range("D2").Value = "'246.58"
# Row 3
$ws.Range("D3").Value = "'21.95"
# Row 4
$ws.Range("D4").Value = "'5.420"
# Row 7
$ws.Range("D7").Value = "'6.320"
# Row 8
$ws.Range("D8").Value = "'0.8072"
# Row 9
$ws.Range("D9").Value = "'0.9542"
$ws.Range("E9").Value = "8FTXTokenFTTBestin24h"
# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1428"
$ws.Range("E10").Value = "9WazirXWRX"
# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07507"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
# Row 12
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03198"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03020"
$ws.Range("E13").Value = "12BitrueCoinBTR"
# Row 14
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "'4.166"
$ws.Range("E14").Value = "13MCDexMCB"
# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09416"
$ws.Range("E15").Value = "14BitMartTokenBMX"
# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001599"
$ws.Range("E16").Value = "15BitForexTokenBF"
# Row 17
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04822"
$ws.Range("E17").Value = "16CoinExTokenCET"
# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005898"
$ws.Range("E18").Value = "17OneONE"
# Row 19
$ws.Range("D19").Value = "'0.006033"
# Row 21
$ws.Range("D21").Value = "'0.0009978"
# Row 24
$ws.Range("D24").Value = "'2.231"
# Row 26
$ws.Range("D26").Value = "'0.1259"
# Row 27
$ws.Range("D27").Value = "'0.0003039"
# Row 40
$ws.Range("D40").Value = "'0.03896"
# Row 41
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1076"
$ws.Range("E41").Value = "40BKEXTokenBKK"
# Row 42
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002660"
$ws.Range("E42").Value = "41CEJICEJI"
# Row 43
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003037"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
# Row 44
$ws.Range("D44").Value = "'0.006705"
# Row 45
$ws.Range("D45").Value = "'0.00005590"
# Row 47
$ws.Range("D47").Value = "'0.3799"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
# Row 48
$ws.Range("D48").Value = "'0.1469"
# Row 49
$ws.Range("D49").Value = "'0.00002099"
